# edit.ps1 - applies the tracked changes described by the diff.
$p = $ppt.ActivePresentation

# --- Slide 2: "Distributed" bullet text update ---
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(1)
$sh2.TextFrame.TextRange.Replace("Use your system, another system and remotes (GitHub). Full history on repositories carried", "Use your system (local repository), another system and remotes (GitHub). Full history on repositories carried in the log")

# --- Slide 7: reposition three shapes ---
$s7 = $p.Slides.Item(7)

$pic8 = $s7.Shapes.Item(5)
$pic8.Top = 274.26220703125

$tb1 = $s7.Shapes.Item(6)
$tb1.Left = 41.98638153076172
$tb1.Top = 363.8707275390625

$tb7 = $s7.Shapes.Item(7)
$tb7.Left = 258.4549865722656
$tb7.Top = 432.86016845703125

# --- Slide 10: renumber steps 12->11, 13->12, 14->13 ---
$s10 = $p.Slides.Item(10)

$sh9 = $s10.Shapes.Item(9)
$sh9.TextFrame.TextRange.Text = "11. Check log to see commits"

$sh10 = $s10.Shapes.Item(10)
$sh10.TextFrame.TextRange.Text = "12. Move back to Master branch, see that the file is ‘as was’"
# Work around an autofit re-layout quirk on this particular shape: editing its
# text causes the runtime to recompute a taller box than PowerPoint originally
# rendered; restore the original (unchanged, per the authoritative diff) height.
$sh10.Height = 24.234411239624023

$sh11 = $s10.Shapes.Item(11)
$sh11.TextFrame.TextRange.Characters(1, 22).Text = "13. Merge Master with "

# --- Slide 11: renumber steps 15->14, 16->15 ---
$s11 = $p.Slides.Item(11)

$sh6 = $s11.Shapes.Item(6)
$sh6.TextFrame.TextRange.Text = "14. Push these changes, from both branches, to the remote repo"

$sh8 = $s11.Shapes.Item(7)
$sh8.TextFrame.TextRange.Text = "15. Check the changes have been made"
